$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row with MathNet.Numerics library licence info
$ws.Range("A10").Value = "MathNet.Numerics"
$ws.Range("C10").Value = "MIT/X11"
$ws.Range("D10").Value = "https://numerics.mathdotnet.com/License.html"

# Thin left/right border around the component + licence cells, like the rest of the table
$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(7).Weight = 2
$ws.Range("A10").Borders.Item(7).ColorIndex = -4105
$ws.Range("A10").Borders.Item(10).LineStyle = 1
$ws.Range("A10").Borders.Item(10).Weight = 2
$ws.Range("A10").Borders.Item(10).ColorIndex = -4105

$ws.Range("C10").Borders.Item(7).LineStyle = 1
$ws.Range("C10").Borders.Item(7).Weight = 2
$ws.Range("C10").Borders.Item(7).ColorIndex = -4105
$ws.Range("C10").Borders.Item(10).LineStyle = 1
$ws.Range("C10").Borders.Item(10).Weight = 2
$ws.Range("C10").Borders.Item(10).ColorIndex = -4105

# Leave the cursor where the author left it before saving
$ws.Range("E16").Select()
